# Sheet1: Ingredients - reorder columns from (name, id, density) to (id, density, name)
$wsIng = $excel.ActiveWorkbook.Worksheets.Item("Ingredients")

$lastRow = $wsIng.UsedRange.Rows.Count

# Capture all existing values (old layout: A=name, B=id, C=density) before overwriting
$names = @()
$ids = @()
$densities = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $names += ,$wsIng.Cells.Item($r, 1).Value()
    $ids += ,$wsIng.Cells.Item($r, 2).Value()
    $densities += ,$wsIng.Cells.Item($r, 3).Value()
}

# New header order: A=id, B=density, C=name
$wsIng.Range("A1").Value = "id"
$wsIng.Range("B1").Value = "density"
$wsIng.Range("C1").Value = "name"

# Rewrite data rows in new column order
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $i + 2
    $wsIng.Cells.Item($r, 1).Value = $ids[$i]
    $wsIng.Cells.Item($r, 2).Value = $densities[$i]
    $wsIng.Cells.Item($r, 3).Value = $names[$i]
}

# Sheet2: Recipes - reorder columns, and add a new "classification" column
$wsRec = $excel.ActiveWorkbook.Worksheets.Item("Recipes")

# Capture old values (old layout: A=name, B=steps, C=id, D=ingredients)
$oldName = $wsRec.Range("A2").Value()
$oldSteps = $wsRec.Range("B2").Value()
$oldId = $wsRec.Range("C2").Value()
$oldIngredients = $wsRec.Range("D2").Value()

# New header order: A=steps, B=ingredients, C=name, D=id, E=classification
$wsRec.Range("A1").Value = "steps"
$wsRec.Range("B1").Value = "ingredients"
$wsRec.Range("C1").Value = "name"
$wsRec.Range("D1").Value = "id"

# Copy header style from an existing styled header cell to the new E1 header
$wsRec.Range("A1").Copy()
$wsRec.Range("E1").PasteSpecial(-4122)
$wsRec.Range("E1").Value = "classification"

# Rewrite data row in new column order
$wsRec.Range("A2").Value = $oldSteps
$wsRec.Range("B2").Value = $oldIngredients
$wsRec.Range("C2").Value = $oldName
$wsRec.Range("D2").Value = $oldId
$wsRec.Range("E2").Value = "1, 1, 1, 1, 1, 1, 1, 0, 1, 0, 0, 0"

Write-Output "done"
